$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: nip becomes a real number (was text), detection_time/created_at refreshed
$ws.Range("A2").Value = 222
$ws.Range("B2").Value = 0.1018543243408203
$ws.Range("C2").Value = "2025-04-22 04:15:25"

# Row 3: new row, numeric nip
$ws.Range("A3").Value = 222
$ws.Range("B3").Value = 0.079498291015625
$ws.Range("C3").Value = "2025-04-22 04:15:34"

# Row 4: new row, nip stored as text "222" (matches original row-2 encoding)
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "222"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").Value = 0.05739951133728027
$ws.Range("C4").Value = "2025-04-22 04:15:41"
